# Apply the "missing Shipment data" commit to the workbook.
#
# Summary of changes:
#  1. Populate the (previously header-only) Shipment sheet with 17 rows of
#     shipment data, including two formula cells and several new lookup
#     strings (carriers / ship options / instructions).
#  2. Populate the remaining quantity column (C) on ContainsPart, plus the
#     vendor/quantity columns (B/C) for the rows that were missing them.
#  3. Move the active sheet / selections around (cosmetic "where the user
#     left the cursor" state) to match the new workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Shipment sheet - fill in the missing rows
# ---------------------------------------------------------------------
$shipment = $wb.Worksheets.Item("Shipment")

# Pre-seed the shared-string table with the new lookup values in the same
# order the source data set introduced them (all carriers, then all ship
# options, then the two instruction blurbs) so the workbook's string table
# comes out in the same order it originally did. These scratch cells get
# cleared at the end of this block; the strings stay alive because the
# real data cells below reference them too.
$scratch = $shipment.Range("Z100:Z108")
$scratch.Cells.Item(1,1).Value = "USPS"
$scratch.Cells.Item(2,1).Value = "FedEx"
$scratch.Cells.Item(3,1).Value = "UPS"
$scratch.Cells.Item(4,1).Value = "CanadaPost"
$scratch.Cells.Item(5,1).Value = "Express"
$scratch.Cells.Item(6,1).Value = "Overnight"
$scratch.Cells.Item(7,1).Value = "Normal"
$scratch.Cells.Item(8,1).Value = "Deliver between 08:00 and 20:00"
$scratch.Cells.Item(9,1).Value = "If noone is home, leave out back"

# columns: A trackingNumber | B orderId | C carrier | D instruction |
#          E shipCost | F shipDate | G shipOption

$shipment.Range("A2").Value = 1000457
$shipment.Range("B2").Value = 1
$shipment.Range("C2").Value = "USPS"
$shipment.Range("D2").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E2").Value = 9.95
$shipment.Range("G2").Value = "Express"

$shipment.Range("A3").Formula = "=A2*2"
$shipment.Range("B3").Value = 2
$shipment.Range("C3").Value = "CanadaPost"
$shipment.Range("D3").Value = "If noone is home, leave out back"
$shipment.Range("E3").Value = 123.45
$shipment.Range("G3").Value = "Normal"

$shipment.Range("A4").Formula = "=A3*2"
$shipment.Range("B4").Value = 3
$shipment.Range("C4").Value = "FedEx"
$shipment.Range("D4").Value = "If noone is home, leave out back"
$shipment.Range("E4").Value = 456.3
$shipment.Range("G4").Value = "Express"

$shipment.Range("A5").Value = 405002
$shipment.Range("B5").Value = 4
$shipment.Range("C5").Value = "USPS"
$shipment.Range("D5").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E5").Value = 32
$shipment.Range("G5").Value = "Normal"

$shipment.Range("A6").Value = 253525
$shipment.Range("B6").Value = 5
$shipment.Range("C6").Value = "CanadaPost"
$shipment.Range("D6").Value = "If noone is home, leave out back"
$shipment.Range("E6").Value = 255.175
$shipment.Range("G6").Value = "Normal"

$shipment.Range("A7").Value = 4653534
$shipment.Range("B7").Value = 6
$shipment.Range("C7").Value = "FedEx"
$shipment.Range("D7").Value = "If noone is home, leave out back"
$shipment.Range("E7").Value = 295.075
$shipment.Range("G7").Value = "Express"

$shipment.Range("A8").Value = 432442
$shipment.Range("B8").Value = 7
$shipment.Range("C8").Value = "USPS"
$shipment.Range("D8").Value = "If noone is home, leave out back"
$shipment.Range("E8").Value = 334.975
$shipment.Range("G8").Value = "Normal"

$shipment.Range("A9").Value = 432423
$shipment.Range("B9").Value = 8
$shipment.Range("C9").Value = "CanadaPost"
$shipment.Range("D9").Value = "If noone is home, leave out back"
$shipment.Range("E9").Value = 374.875
$shipment.Range("G9").Value = "Overnight"

$shipment.Range("A10").Value = 35435356
$shipment.Range("B10").Value = 9
$shipment.Range("C10").Value = "FedEx"
$shipment.Range("D10").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E10").Value = 414.775
$shipment.Range("G10").Value = "Overnight"

$shipment.Range("A11").Value = 432465
$shipment.Range("B11").Value = 10
$shipment.Range("C11").Value = "USPS"
$shipment.Range("D11").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E11").Value = 454.675
$shipment.Range("G11").Value = "Normal"

$shipment.Range("A12").Value = 353656
$shipment.Range("B12").Value = 11
$shipment.Range("C12").Value = "CanadaPost"
$shipment.Range("D12").Value = "If noone is home, leave out back"
$shipment.Range("E12").Value = 494.575
$shipment.Range("G12").Value = "Express"

$shipment.Range("A13").Value = 854
$shipment.Range("B13").Value = 12
$shipment.Range("C13").Value = "UPS"
$shipment.Range("D13").Value = "If noone is home, leave out back"
$shipment.Range("E13").Value = 534.475
$shipment.Range("G13").Value = "Overnight"

$shipment.Range("A14").Value = 12325
$shipment.Range("B14").Value = 13
$shipment.Range("C14").Value = "USPS"
$shipment.Range("D14").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E14").Value = 574.375
$shipment.Range("G14").Value = "Normal"

$shipment.Range("A15").Value = 533
$shipment.Range("B15").Value = 14
$shipment.Range("C15").Value = "CanadaPost"
$shipment.Range("D15").Value = "If noone is home, leave out back"
$shipment.Range("E15").Value = 614.275
$shipment.Range("G15").Value = "Express"

$shipment.Range("A16").Value = 53
$shipment.Range("B16").Value = 15
$shipment.Range("C16").Value = "UPS"
$shipment.Range("D16").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E16").Value = 654.175
$shipment.Range("G16").Value = "Normal"

$shipment.Range("A17").Value = 525
$shipment.Range("B17").Value = 16
$shipment.Range("C17").Value = "USPS"
$shipment.Range("D17").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E17").Value = 694.075
$shipment.Range("G17").Value = "Overnight"

$shipment.Range("A18").Value = 25252
$shipment.Range("B18").Value = 17
$shipment.Range("C18").Value = "USPS"
$shipment.Range("D18").Value = "Deliver between 08:00 and 20:00"
$shipment.Range("E18").Value = 733.975
$shipment.Range("G18").Value = "Overnight"

# Drop the scratch cells now that every new string is referenced by real data.
$scratch.ClearContents()

# ---------------------------------------------------------------------
# 2. ContainsPart sheet - fill in the quantity column + missing rows
# ---------------------------------------------------------------------
$containsPart = $wb.Worksheets.Item("ContainsPart")

$containsPart.Range("C2").Value = 1
$containsPart.Range("C3").Value = 1
$containsPart.Range("C4").Value = 1
$containsPart.Range("C5").Value = 1
$containsPart.Range("C6").Value = 1
$containsPart.Range("C7").Value = 1
$containsPart.Range("C8").Value = 1

$containsPart.Range("B9").Value = 1
$containsPart.Range("C9").Value = 2

$containsPart.Range("B10").Value = 2
$containsPart.Range("C10").Value = 3

$containsPart.Range("B11").Value = 3
$containsPart.Range("C11").Value = 3

$containsPart.Range("B12").Value = 4
$containsPart.Range("C12").Value = 2

$containsPart.Range("B13").Value = 5
$containsPart.Range("C13").Value = 1

$containsPart.Range("B14").Value = 3
$containsPart.Range("C14").Value = 1

$containsPart.Range("B15").Value = 2
$containsPart.Range("C15").Value = 1

$containsPart.Range("B16").Value = 1
$containsPart.Range("C16").Value = 1

$containsPart.Range("B17").Value = 1
$containsPart.Range("C17").Value = 1

$containsPart.Range("B18").Value = 1
$containsPart.Range("C18").Value = 1

# ---------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping
# ---------------------------------------------------------------------
$wb.Worksheets.Item("RatesVendor").Range("B29").Select()
$wb.Worksheets.Item("PartOrder").Range("D12").Select()
$wb.Worksheets.Item("ContainsPart").Range("C19").Select()
$wb.Worksheets.Item("Part").Range("D17").Select()

# Shipment becomes the active tab, selection moves to D20.
$shipment.Activate()
$shipment.Range("D20").Select()
